# Auto-generated edit script applying the Adamantoise_Profits data refresh
# (scheduled market-price runner) diff to before.xlsx
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 910283
$ws.Range("I6").Value = 1112489.4
$ws.Range("K6").Value = 3337468.2
$ws.Range("M6").Value = -3337356.2
$ws.Range("H19").Value = 1197.5454
$ws.Range("I19").Value = 1458.1428
$ws.Range("J19").Value = 741.5
$ws.Range("K19").Value = 1458.1428
$ws.Range("L19").Value = 741.5
$ws.Range("M19").Value = -1283.1428
$ws.Range("N19").Value = -1091.5
$ws.Range("H132").Value = 2297.639
$ws.Range("I132").Value = 1920.44
$ws.Range("K132").Value = 5761.32
$ws.Range("M132").Value = -3231.32
$ws.Range("H138").Value = 2537.41
$ws.Range("I138").Value = 1161.75
$ws.Range("J138").Value = 2725
$ws.Range("K138").Value = 3485.25
$ws.Range("L138").Value = 8175
$ws.Range("M138").Value = 1654.75
$ws.Range("N138").Value = -18455

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 60961.332
$ws.Range("J52").Value = 60961.332
$ws.Range("L52").Value = 60961.332
$ws.Range("N52").Value = -61597.332
$ws.Range("H88").Value = 1771.2
$ws.Range("I88").Value = 2150
$ws.Range("J88").Value = 1203
$ws.Range("K88").Value = 2150
$ws.Range("L88").Value = 1203
$ws.Range("M88").Value = -1744
$ws.Range("N88").Value = -2015
$ws.Range("H91").Value = 1771.2
$ws.Range("I91").Value = 2150
$ws.Range("J91").Value = 1203
$ws.Range("K91").Value = 2150
$ws.Range("L91").Value = 1203
$ws.Range("M91").Value = -746
$ws.Range("N91").Value = -4011
$ws.Range("H97").Value = 1288.8379
$ws.Range("J97").Value = 1599.4286
$ws.Range("L97").Value = 1599.4286
$ws.Range("N97").Value = -2591.4286
$ws.Range("H133").Value = 62249.75
$ws.Range("J133").Value = 62249.75
$ws.Range("L133").Value = 62249.75
$ws.Range("N133").Value = -67309.75

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 528.6
$ws.Range("J80").Value = 597.7143
$ws.Range("L80").Value = 597.7143
$ws.Range("N80").Value = -2593.7143
$ws.Range("H83").Value = 528.6
$ws.Range("J83").Value = 597.7143
$ws.Range("L83").Value = 2988.5715
$ws.Range("N83").Value = -12972.5715
$ws.Range("H99").Value = 2164.8235
$ws.Range("J99").Value = 3200
$ws.Range("L99").Value = 3200
$ws.Range("N99").Value = -6196

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2778.8518
$ws.Range("I58").Value = 2563.7083
$ws.Range("J58").Value = 4500
$ws.Range("K58").Value = 2563.7083
$ws.Range("L58").Value = 4500
$ws.Range("M58").Value = -2360.7083
$ws.Range("N58").Value = -4906
$ws.Range("H62").Value = 2129.5715
$ws.Range("I62").Value = 2182
$ws.Range("J62").Value = 1815
$ws.Range("K62").Value = 2182
$ws.Range("L62").Value = 1815
$ws.Range("M62").Value = -1558
$ws.Range("N62").Value = -3063
$ws.Range("H65").Value = 2129.5715
$ws.Range("I65").Value = 2182
$ws.Range("J65").Value = 1815
$ws.Range("K65").Value = 10910
$ws.Range("L65").Value = 9075
$ws.Range("M65").Value = -7790
$ws.Range("N65").Value = -15315
$ws.Range("H105").Value = 2082.2778
$ws.Range("I105").Value = 2069.1538
$ws.Range("J105").Value = 2116.4
$ws.Range("K105").Value = 2069.1538
$ws.Range("L105").Value = 2116.4
$ws.Range("M105").Value = -322.1538
$ws.Range("N105").Value = -5610.4
$ws.Range("H136").Value = 2778.8518
$ws.Range("I136").Value = 2563.7083
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 7691.124899999999
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -5141.124899999999
$ws.Range("N136").Value = -18600

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 573.9231
$ws.Range("I44").Value = 469.5
$ws.Range("J44").Value = 741
$ws.Range("K44").Value = 1408.5
$ws.Range("L44").Value = 2223
$ws.Range("M44").Value = -1010.5
$ws.Range("N44").Value = -3019
$ws.Range("H47").Value = 45501000
$ws.Range("I47").Value = 91000000
$ws.Range("K47").Value = 273000000
$ws.Range("M47").Value = -272999569
$ws.Range("H58").Value = 2249.25
$ws.Range("I58").Value = 2249.25
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 6747.75
$ws.Range("L58").Value = 0
$ws.Range("N58").Value = -6619.75
$ws.Range("H60").Value = 522.3333
$ws.Range("I60").Value = 424
$ws.Range("J60").Value = 601
$ws.Range("K60").Value = 1272
$ws.Range("L60").Value = 1803
$ws.Range("M60").Value = -1021
$ws.Range("N60").Value = -2305
$ws.Range("H131").Value = 1620.5227
$ws.Range("I131").Value = 1043.9166
$ws.Range("K131").Value = 3131.7498
$ws.Range("M131").Value = 1908.2502

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3189.125
$ws.Range("I80").Value = 3201
$ws.Range("J80").Value = 3177.25
$ws.Range("K80").Value = 3201
$ws.Range("L80").Value = 3177.25
$ws.Range("M80").Value = -2203
$ws.Range("N80").Value = -5173.25
$ws.Range("H83").Value = 3189.125
$ws.Range("I83").Value = 3201
$ws.Range("J83").Value = 3177.25
$ws.Range("K83").Value = 16005
$ws.Range("L83").Value = 15886.25
$ws.Range("M83").Value = -11013
$ws.Range("N83").Value = -25870.25
$ws.Range("H122").Value = 1597.0625
$ws.Range("I122").Value = 1553.8462
$ws.Range("J122").Value = 1784.3334
$ws.Range("K122").Value = 4661.5386
$ws.Range("L122").Value = 5353.0002
$ws.Range("M122").Value = -2211.5386
$ws.Range("N122").Value = -10253.0002
$ws.Range("H126").Value = 2902.3635
$ws.Range("I126").Value = 2813.5
$ws.Range("J126").Value = 3057.875
$ws.Range("K126").Value = 8440.5
$ws.Range("L126").Value = 9173.625
$ws.Range("M126").Value = -5970.5
$ws.Range("N126").Value = -14113.625

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1604.5714
$ws.Range("I68").Value = 1677.75
$ws.Range("K68").Value = 1677.75
$ws.Range("M68").Value = -928.75
$ws.Range("H71").Value = 1604.5714
$ws.Range("I71").Value = 1677.75
$ws.Range("K71").Value = 8388.75
$ws.Range("M71").Value = -4644.75
$ws.Range("H93").Value = 1603.5416
$ws.Range("I93").Value = 1318.6
$ws.Range("J93").Value = 2078.4443
$ws.Range("K93").Value = 1318.6
$ws.Range("L93").Value = 2078.4443
$ws.Range("M93").Value = -70.59999999999991
$ws.Range("N93").Value = -4574.4443

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 59478.555
$ws.Range("J81").Value = 5220.8887
$ws.Range("L81").Value = 10441.7774
$ws.Range("N81").Value = -12563.7774
$ws.Range("H84").Value = 59478.555
$ws.Range("J84").Value = 5220.8887
$ws.Range("M84").Value = -62816.887
$ws.Range("N84").Value = -62816.887

# ---- Structural fixups ----
# CUL row 58: LeveProfitNQ (M58) cell is removed entirely in the refreshed
# data (currentAveragePriceNQ now equals the HQ average, so the NQ-profit
# column is dropped) while N58 (LeveProfitHQ) keeps its new value.
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M58").ClearContents()
